$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the rows that will no longer be used (A9:A24)
$ws.Range("A9:A24").ClearContents()

# Consolidate each card's attributes into a single Python-tuple-style string per row
$ws.Range("A2").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A3").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A4").Value = "('Lightning Bolt', ['{R}', 'Instant', 'Lightning Bolt deals 3 damage to any target.'])"
$ws.Range("A5").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A6").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A7").Value = "('Sol Ring', ['{1}', 'Artifact', '{T}: Add {C}{C}.'])"
$ws.Range("A8").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"
